# Apply cryptos.xlsx data refresh: update prices/percentages, and
# swap the Filecoin/LEO rows (28 and 29) per the commit's data update.
# Values that look like plain numbers (e.g. "1.00", "138.42") are
# prefixed with a leading apostrophe so Excel stores them as literal
# text instead of auto-converting them to numeric values (matching
# the source data, which keeps these as text strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.088.41"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.779.42"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'430.22"
$ws.Range("E5").Value = "  +5.37%  "
$ws.Range("D6").Value = "'138.42"
$ws.Range("E6").Value = "  +4.77%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.734"
$ws.Range("D10").Value = "'0.151"
$ws.Range("E10").Value = "  -8.72%  "
$ws.Range("D11").Value = "'0.0000310"
$ws.Range("E11").Value = "  -12.81%  "
$ws.Range("D12").Value = "'42.52"
$ws.Range("E12").Value = "  +3.57%  "
$ws.Range("D13").Value = "'10.41"
$ws.Range("E13").Value = "  +4.80%  "
$ws.Range("D14").Value = "4.361.58"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "'14.94"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "3.770.43"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'19.94"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("E19").Value = "  +5.82%  "
$ws.Range("D20").Value = "66.145.41"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "'404.75"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("E22").Value = "  +4.07%  "
$ws.Range("D23").Value = "'3.27"
$ws.Range("E23").Value = "  +6.70%  "
$ws.Range("D24").Value = "'84.64"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").Value = "'10.01"
$ws.Range("E25").Value = "  +35.45%  "
$ws.Range("D26").Value = "'36.54"
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("E27").Value = "  +5.73%  "

# Row 28 and 29 swap: Filecoin moves to row 28, LEO moves to row 29
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'9.88"
$ws.Range("E28").Value = "  +5.89%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "'5.52"
$ws.Range("E29").Value = "  -3.30%  "

$ws.Range("D30").Value = "'0.137"
$ws.Range("E30").Value = "  +13.71%  "
$ws.Range("D31").Value = "'13.78"
$ws.Range("E31").Value = "  +11.86%  "
$ws.Range("D32").Value = "'704.10"
$ws.Range("E32").Value = "  -3.78%  "
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").Value = "'41.33"
$ws.Range("E34").Value = "  +6.27%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'5.75"
$ws.Range("E36").Value = "  +35.30%  "
$ws.Range("D37").Value = "'0.148"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("D39").Value = "'0.0471"
$ws.Range("E39").Value = "  +2.60%  "
$ws.Range("D40").Value = "'2.80"
$ws.Range("E40").Value = "  +40.32%  "
$ws.Range("D41").Value = "'2.99"
$ws.Range("E41").Value = "  +5.41%  "
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "0.0₃0670"
$ws.Range("E44").Value = "  -8.97%  "
$ws.Range("D45").Value = "'0.329"
$ws.Range("E45").Value = "  +11.57%  "
$ws.Range("D46").Value = "'3.20"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("E47").Value = "  +2.43%  "
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "'139.35"
$ws.Range("E50").Value = "  -4.54%  "
$ws.Range("E51").Value = "  -0.29%  "
